$wb = $excel.ActiveWorkbook

# --- Sheets ---
$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 and de-de!H2 share the same underlying string value.
$overview.Range("G2").Value = "2016-10-21 01:04:38"
$dede.Range("H2").Value     = "2016-10-21 01:04:38"

# zh-cn!H2 gets its own updated timestamp.
$zhcn.Range("H2").Value = "2016-10-21 01:04:27"

# --- Column widths (zh-cn / de-de "Status" column widened; Overview's matching columns too) ---
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth     = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth     = 16.333333333333332
